$d = $word.ActiveDocument

# Update the date heading
$dateRange = $d.Paragraphs(1).Range
$dateRange.MoveEnd(1, -1) | Out-Null
$dateRange.Text = "2026-01-26 Monday"

# Update table cell values (row-major order)
$newValues = @(
    "28+25=",
    "76-27=",
    "12+59=",
    "85-47=",
    "38+3=",
    "92-76=",
    "90-51=",
    "14+48=",
    "74-46=",
    "55+38=",
    "46+17=",
    "47-28=",
    "15+38=",
    "82-14=",
    "91-15=",
    "62-34=",
    "35+27=",
    "48+34=",
    "64+17=",
    "38+47=",
    "37+55=",
    "73-35=",
    "19+38=",
    "43-6=",
    "71-36=",
    "52-49=",
    "16+18=",
    "54-29=",
    "36-19=",
    "18+76=",
    "60-38=",
    "29+62=",
    "67+19=",
    "36+15=",
    "80-47=",
    "64-58=",
    "96-7=",
    "15+68=",
    "88-9=",
    "92-13=",
    "73-48=",
    "54-38=",
    "35+8=",
    "30-27=",
    "63-58=",
    "25+7=",
    "47+35=",
    "20-15=",
    "64-27=",
    "58+18=",
    "68+19=",
    "54+29=",
    "28+25=",
    "19+33=",
    "19+15=",
    "61-17=",
    "9+2=",
    "84-37=",
    "25+19=",
    "91-46=",
    "14+29=",
    "82-17=",
    "75-27=",
    "80-23=",
    "81-24=",
    "37+25=",
    "61-44=",
    "70-11=",
    "68+24=",
    "88-29=",
    "91-2=",
    "86-38=",
    "29+6=",
    "6+16=",
    "50-17=",
    "66+8=",
    "56-47=",
    "2+89=",
    "42-15=",
    "43-34=",
    "46+45=",
    "18+5=",
    "41-38=",
    "64-5=",
    "64-7=",
    "61-46=",
    "94-56=",
    "3+29=",
    "68+16=",
    "38+54=",
    "19+72=",
    "80-58=",
    "88+4=",
    "29+65=",
    "42-14=",
    "93-68=",
    "41-5=",
    "65-8=",
    "71-58=",
    "81-59="
)

$t = $d.Tables(1)
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cr = $cell.Range
        $cr.MoveEnd(1, -2) | Out-Null
        $cr.Text = $newValues[$idx]
        $idx++
    }
}
